$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("FEB-22")

$ws.Cells.Item(9, 1).Value = 5
$ws.Cells.Item(9, 2).Value = 44596
$ws.Cells.Item(9, 3).Value = "RPA GSS"
$ws.Cells.Item(9, 4).Value = "1. DRS Monthly task has been tested for the 5000 records and it is running smoothly"
$ws.Cells.Item(9, 5).Value = 1
$ws.Cells.Item(9, 6).Value = "Completed"

$ws.Cells.Item(10, 4).Value = "2. Implementation of public holidays has been done Warranty Daily task, tested  and running smoothly"
$ws.Cells.Item(10, 5).Value = 1
$ws.Cells.Item(10, 6).Value = "Completed"

$ws.Cells.Item(11, 4).Value = "3. Implementation of public holidays has been done Activity Daily task, tested  and running smoothly"
$ws.Cells.Item(11, 5).Value = 1
$ws.Cells.Item(11, 6).Value = "Completed"

$ws.Cells.Item(12, 4).Value = "4. Implementation of public holidays at GRD task(GRSummary with Details task) is work in progress"
$ws.Cells.Item(12, 5).Value = 0.1
$ws.Cells.Item(12, 6).Value = "WIP"

$ws.Cells.Item(13, 1).Value = 6
$ws.Cells.Item(13, 2).Value = 44596
$ws.Cells.Item(13, 3).Value = "RPA RLOGIC"
$ws.Cells.Item(13, 4).Value = "5. Cutomization at P&L is work in progress  "
$ws.Cells.Item(13, 5).Value = 0.3
$ws.Cells.Item(13, 6).Value = "WIP"

$ws.Range("D19").Select()
